# Auto-generated Excel COM-interop script
# Applies updated market price / profit values to the leve profit sheets
# as produced by the scheduled data-refresh runner.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (49 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 60607.75
$ws.Range("J17").Value = 62325.113
$ws.Range("L17").Value = 186975.339
$ws.Range("N17").Value = -187311.339
$ws.Range("H62").Value = 3560.8333
$ws.Range("I62").Value = 2900.0
$ws.Range("J62").Value = 3891.25
$ws.Range("K62").Value = 2900.0
$ws.Range("L62").Value = 3891.25
$ws.Range("M62").Value = -2276.0
$ws.Range("N62").Value = -5139.25
$ws.Range("H65").Value = 3560.8333
$ws.Range("I65").Value = 2900.0
$ws.Range("J65").Value = 3891.25
$ws.Range("K65").Value = 14500.0
$ws.Range("L65").Value = 19456.25
$ws.Range("M65").Value = -11380.0
$ws.Range("N65").Value = -25696.25
$ws.Range("H82").Value = 4677.0
$ws.Range("H85").Value = 4677.0
$ws.Range("H111").Value = 568.8333
$ws.Range("I111").Value = 543.6667
$ws.Range("J111").Value = 644.3333
$ws.Range("K111").Value = 1631.0001
$ws.Range("L111").Value = 1932.9999
$ws.Range("M111").Value = 1435.9999
$ws.Range("N111").Value = -8066.9999
$ws.Range("H116").Value = 5248.3
$ws.Range("I116").Value = 5204.4443
$ws.Range("K116").Value = 5204.4443
$ws.Range("M116").Value = -1762.4443
$ws.Range("H125").Value = 1563.6364
$ws.Range("I125").Value = 1600.0
$ws.Range("K125").Value = 14400.0
$ws.Range("M125").Value = -11940.0
$ws.Range("H135").Value = 1467.4
$ws.Range("I135").Value = 1557.6923
$ws.Range("J135").Value = 1299.7142
$ws.Range("K135").Value = 14019.2307
$ws.Range("L135").Value = 11697.4278
$ws.Range("M135").Value = -11484.2307
$ws.Range("N135").Value = -16767.4278
$ws.Range("H141").Value = 804423.1
$ws.Range("I141").Value = 3871.375
$ws.Range("J141").Value = 1160223.9
$ws.Range("K141").Value = 11614.125
$ws.Range("L141").Value = 3480671.7
$ws.Range("M141").Value = -6434.125
$ws.Range("N141").Value = -3491031.7

# --- Sheet: ARM (33 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 27780956.0
$ws.Range("I2").Value = 50002320.0
$ws.Range("J2").Value = 4249.75
$ws.Range("K2").Value = 50002320.0
$ws.Range("L2").Value = 4249.75
$ws.Range("M2").Value = -50002207.0
$ws.Range("N2").Value = -4475.75
$ws.Range("H63").Value = 2409.762
$ws.Range("I63").Value = 1631.625
$ws.Range("K63").Value = 1631.625
$ws.Range("M63").Value = -945.625
$ws.Range("H66").Value = 2409.762
$ws.Range("I66").Value = 1631.625
$ws.Range("K66").Value = 8158.125
$ws.Range("M66").Value = -4726.125
$ws.Range("H97").Value = 458.41177
$ws.Range("I97").Value = 365.7586
$ws.Range("K97").Value = 365.7586
$ws.Range("M97").Value = 130.2414
$ws.Range("H116").Value = 27780956.0
$ws.Range("I116").Value = 50002320.0
$ws.Range("J116").Value = 4249.75
$ws.Range("K116").Value = 50002320.0
$ws.Range("L116").Value = 4249.75
$ws.Range("M116").Value = -50000026.0
$ws.Range("N116").Value = -8837.75
$ws.Range("H132").Value = 2072.7058
$ws.Range("I132").Value = 1780.439
$ws.Range("J132").Value = 3271.0
$ws.Range("K132").Value = 5341.317
$ws.Range("L132").Value = 9813.0
$ws.Range("M132").Value = -2811.317
$ws.Range("N132").Value = -14873.0

# --- Sheet: BSM (35 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 27780956.0
$ws.Range("I3").Value = 50002320.0
$ws.Range("J3").Value = 4249.75
$ws.Range("K3").Value = 50002320.0
$ws.Range("L3").Value = 4249.75
$ws.Range("M3").Value = -50002206.0
$ws.Range("N3").Value = -4477.75
$ws.Range("H22").Value = 217.28572
$ws.Range("I22").Value = 228.5
$ws.Range("J22").Value = 150.0
$ws.Range("K22").Value = 228.5
$ws.Range("L22").Value = 150.0
$ws.Range("M22").Value = -55.5
$ws.Range("N22").Value = -496.0
$ws.Range("H99").Value = 2172.7222
$ws.Range("I99").Value = 1473.9333
$ws.Range("J99").Value = 5666.6665
$ws.Range("K99").Value = 1473.9333
$ws.Range("L99").Value = 5666.6665
$ws.Range("M99").Value = 24.06670000000008
$ws.Range("N99").Value = -8662.6665
$ws.Range("H105").Value = 2950.0
$ws.Range("I105").Value = 3000.0
$ws.Range("J105").Value = 2900.0
$ws.Range("K105").Value = 3000.0
$ws.Range("L105").Value = 2900.0
$ws.Range("M105").Value = -1253.0
$ws.Range("N105").Value = -6394.0
$ws.Range("H134").Value = 2368.0645
$ws.Range("I134").Value = 1285.6111
$ws.Range("J134").Value = 3866.8462
$ws.Range("K134").Value = 3856.8333
$ws.Range("L134").Value = 11600.5386
$ws.Range("M134").Value = -1321.8333
$ws.Range("N134").Value = -16670.5386

# --- Sheet: CRP (50 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1757262.1
$ws.Range("J31").Value = 7230.6924
$ws.Range("L31").Value = 7230.6924
$ws.Range("N31").Value = -7820.6924
$ws.Range("H34").Value = 1757262.1
$ws.Range("J34").Value = 7230.6924
$ws.Range("L34").Value = 7230.6924
$ws.Range("N34").Value = -7634.6924
$ws.Range("H58").Value = 10419835.0
$ws.Range("I58").Value = 1817.3793
$ws.Range("J58").Value = 26321020.0
$ws.Range("K58").Value = 1817.3793
$ws.Range("L58").Value = 26321020.0
$ws.Range("M58").Value = -1614.3793
$ws.Range("N58").Value = -26321426.0
$ws.Range("H62").Value = 4585.0
$ws.Range("I62").Value = 3577.5
$ws.Range("J62").Value = 6600.0
$ws.Range("K62").Value = 3577.5
$ws.Range("L62").Value = 6600.0
$ws.Range("M62").Value = -2953.5
$ws.Range("N62").Value = -7848.0
$ws.Range("H65").Value = 4585.0
$ws.Range("I65").Value = 3577.5
$ws.Range("J65").Value = 6600.0
$ws.Range("K65").Value = 17887.5
$ws.Range("L65").Value = 33000.0
$ws.Range("M65").Value = -14767.5
$ws.Range("N65").Value = -39240.0
$ws.Range("H122").Value = 1440.675
$ws.Range("I122").Value = 1197.7188
$ws.Range("J122").Value = 2412.5
$ws.Range("K122").Value = 3593.1564
$ws.Range("L122").Value = 7237.5
$ws.Range("M122").Value = -1143.1564
$ws.Range("N122").Value = -12137.5
$ws.Range("H132").Value = 2006.6976
$ws.Range("I132").Value = 1926.4286
$ws.Range("J132").Value = 2083.318
$ws.Range("K132").Value = 5779.2858
$ws.Range("L132").Value = 6249.954000000001
$ws.Range("M132").Value = -3249.2858
$ws.Range("N132").Value = -11309.954
$ws.Range("H136").Value = 10419835.0
$ws.Range("I136").Value = 1817.3793
$ws.Range("J136").Value = 26321020.0
$ws.Range("K136").Value = 5452.1379
$ws.Range("L136").Value = 78963060.0
$ws.Range("M136").Value = -2902.1379
$ws.Range("N136").Value = -78968160.0

# --- Sheet: CUL (4 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2550.4
$ws.Range("J131").Value = 2680.5557
$ws.Range("L131").Value = 8041.6671
$ws.Range("N131").Value = -18121.6671

# --- Sheet: GSM (11 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 31543.234
$ws.Range("I102").Value = 1809.55
$ws.Range("J102").Value = 74019.93
$ws.Range("K102").Value = 1809.55
$ws.Range("L102").Value = 74019.93
$ws.Range("M102").Value = -187.55
$ws.Range("N102").Value = -77263.93
$ws.Range("H122").Value = 3986.6785
$ws.Range("I122").Value = 2654.5293
$ws.Range("K122").Value = 7963.5879
$ws.Range("M122").Value = -5513.5879

# --- Sheet: LTW (35 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1925.5264
$ws.Range("I7").Value = 1212.6666
$ws.Range("J7").Value = 2254.5386
$ws.Range("K7").Value = 1212.6666
$ws.Range("L7").Value = 2254.5386
$ws.Range("M7").Value = -1100.6666
$ws.Range("N7").Value = -2478.5386
$ws.Range("H40").Value = 8385.714
$ws.Range("I40").Value = 11175.0
$ws.Range("J40").Value = 4666.6665
$ws.Range("K40").Value = 11175.0
$ws.Range("L40").Value = 4666.6665
$ws.Range("M40").Value = -11039.0
$ws.Range("N40").Value = -4938.6665
$ws.Range("H61").Value = 1179.8125
$ws.Range("I61").Value = 529.0
$ws.Range("J61").Value = 4000.0
$ws.Range("K61").Value = 529.0
$ws.Range("L61").Value = 4000.0
$ws.Range("M61").Value = -327.0
$ws.Range("N61").Value = -4404.0
$ws.Range("H113").Value = 1179.8125
$ws.Range("I113").Value = 529.0
$ws.Range("J113").Value = 4000.0
$ws.Range("K113").Value = 529.0
$ws.Range("L113").Value = 4000.0
$ws.Range("M113").Value = 1641.0
$ws.Range("N113").Value = -8340.0
$ws.Range("H126").Value = 1925.5264
$ws.Range("I126").Value = 1212.6666
$ws.Range("J126").Value = 2254.5386
$ws.Range("K126").Value = 3637.9998
$ws.Range("L126").Value = 6763.6158
$ws.Range("M126").Value = -1167.9998
$ws.Range("N126").Value = -11703.6158

# --- Sheet: WVR (32 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1715.5
$ws.Range("I107").Value = 384.8
$ws.Range("J107").Value = 3933.3333
$ws.Range("K107").Value = 1154.4
$ws.Range("L107").Value = 11799.9999
$ws.Range("M107").Value = 765.5999999999999
$ws.Range("N107").Value = -15639.9999
$ws.Range("H113").Value = 1404.8077
$ws.Range("I113").Value = 809.6923
$ws.Range("J113").Value = 1999.9231
$ws.Range("K113").Value = 2429.0769
$ws.Range("L113").Value = 5999.7693
$ws.Range("M113").Value = -259.0769
$ws.Range("N113").Value = -10339.7693
$ws.Range("H122").Value = 590594.8
$ws.Range("I122").Value = 835235.5
$ws.Range("J122").Value = 3457.2
$ws.Range("K122").Value = 2505706.5
$ws.Range("L122").Value = 10371.6
$ws.Range("M122").Value = -2503256.5
$ws.Range("N122").Value = -15271.6
$ws.Range("H126").Value = 3705717.5
$ws.Range("I126").Value = 1278.4375
$ws.Range("J126").Value = 9093993.0
$ws.Range("K126").Value = 3835.3125
$ws.Range("L126").Value = 27281979.0
$ws.Range("M126").Value = -1365.3125
$ws.Range("N126").Value = -27286919.0
$ws.Range("H132").Value = 242675.72
$ws.Range("I132").Value = 401296.72
$ws.Range("K132").Value = 1203890.16
$ws.Range("M132").Value = -1201360.16
